$d = $word.ActiveDocument
$d.Content.Find.Execute("Bob Weisberg", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Robert H. Weisberg", 2)
